$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the trailing empty placeholder rows (13, then 12) - row 11 stays and
# gets new content below.
$ws.Rows("13").Delete()
$ws.Rows("12").Delete()

# Insert two new columns at D:E (old column D -> F, old column J -> L)
$ws.Range("D1:E1").EntireColumn.Insert()

# New header cells
$ws.Range("D1").Value = "Type"
$ws.Range("E1").Value = " URL"

# Data rows: Type / URL columns for the fall2017 project rows
$ws.Range("D4").Value = "Web"
$ws.Range("E4").Value = "D:\ResearchSpace\Repositories\577 projects\fall2017\projects\f17team01\Valuation"

$ws.Range("D5").Value = "Web"
$ws.Range("E5").Value = "D:\ResearchSpace\Repositories\577 projects\fall2017\projects\f17team02\Foundations\DCP"

$ws.Range("D6").Value = "Web"
$ws.Range("E6").Value = "D:\ResearchSpace\Repositories\577 projects\fall2017\projects\f17team03\Foundations"

$ws.Range("D7").Value = "App"

$ws.Range("D8").Value = "App"

$ws.Range("D9").Value = "App"
$ws.Range("E9").Value = "D:\ResearchSpace\Repositories\577 projects\fall2017\projects\f17team06\Valuation\TA"

$ws.Range("D10").Value = "Web"

$ws.Range("D11").Value = 2018

# Drop the stale cached row heights so rows fall back to the sheet default
# (matches the re-saved file, which no longer pins a per-row height)
$ws.Rows("1:11").AutoFit() | Out-Null

# Match the selection shown in the edited file
$ws.Range("D11").Select()
